$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 20
$ws.Range("H20").Value = 297.75
$ws.Range("I20").Value = 297.75
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 297.75
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
# Row 33
$ws.Range("H33").Value = 113.25
$ws.Range("I33").Value = 84.5
$ws.Range("J33").Value = 199.5
$ws.Range("K33").Value = 84.5
$ws.Range("L33").Value = 199.5
$ws.Range("M33").Value = 144.5
$ws.Range("N33").Value = -657.5
# Row 35
$ws.Range("H35").Value = 297.75
$ws.Range("I35").Value = 297.75
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 297.75
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
# Row 41
$ws.Range("H41").Value = 535
$ws.Range("I41").Value = 609
$ws.Range("J41").Value = 411.66666
$ws.Range("K41").Value = 609
$ws.Range("L41").Value = 411.66666
$ws.Range("M41").Value = -169
$ws.Range("N41").Value = -1291.66666
# Row 46
$ws.Range("H46").Value = 17
$ws.Range("I46").Value = 17
$ws.Range("K46").Value = 51
$ws.Range("M46").Value = 68
# Row 49
$ws.Range("H49").Value = 299.5
$ws.Range("I49").Value = 299.5
$ws.Range("K49").Value = 898.5
$ws.Range("M49").Value = -762.5
# Row 60
$ws.Range("H60").Value = 17
$ws.Range("I60").Value = 17
$ws.Range("K60").Value = 51
$ws.Range("M60").Value = 433
# Row 76
$ws.Range("H76").Value = 4150
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
# Row 79
$ws.Range("H79").Value = 4150
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
# Row 92
$ws.Range("H92").Value = 153.0625
$ws.Range("I92").Value = 146.35715
$ws.Range("K92").Value = 146.35715
$ws.Range("M92").Value = 1101.64285
# Row 135
$ws.Range("H135").Value = 741.3214
$ws.Range("I135").Value = 658.4761999999999
$ws.Range("K135").Value = 5926.2858
$ws.Range("M135").Value = -3391.2858
# Row 137
$ws.Range("H137").Value = 6292.533
$ws.Range("I137").Value = 2341.4285
$ws.Range("K137").Value = 7024.2855
$ws.Range("M137").Value = -4474.2855

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3032754.5
$ws.Range("I32").Value = 567.75
$ws.Range("K32").Value = 567.75
$ws.Range("M32").Value = -280.75
# Row 61
$ws.Range("H61").Value = 2041
$ws.Range("I61").Value = 1931.9231
$ws.Range("K61").Value = 1931.9231
$ws.Range("M61").Value = -1719.9231
# Row 74
$ws.Range("H74").Value = 3866.7222
$ws.Range("I74").Value = 3594.125
$ws.Range("K74").Value = 3594.125
$ws.Range("M74").Value = -2720.125
# Row 77
$ws.Range("H77").Value = 3866.7222
$ws.Range("I77").Value = 3594.125
$ws.Range("K77").Value = 17970.625
$ws.Range("M77").Value = -13602.625
# Row 105
$ws.Range("H105").Value = 24995
$ws.Range("J105").Value = 24995
$ws.Range("L105").Value = 24995
$ws.Range("N105").Value = -31983
# Row 110
$ws.Range("H110").Value = 100003100
$ws.Range("I110").Value = 333337000
$ws.Range("J110").Value = 2857.1428
$ws.Range("K110").Value = 333337000
$ws.Range("L110").Value = 2857.1428
$ws.Range("M110").Value = -333334955
$ws.Range("N110").Value = -6947.1428
# Row 112
$ws.Range("H112").Value = 17024.666
$ws.Range("J112").Value = 17024.666
$ws.Range("L112").Value = 17024.666
$ws.Range("N112").Value = -19978.666
# Row 122
$ws.Range("H122").Value = 2009.9286
$ws.Range("I122").Value = 1463.4546
$ws.Range("K122").Value = 4390.3638
$ws.Range("M122").Value = -1940.3638
# Row 132
$ws.Range("H132").Value = 4132.4546
$ws.Range("I132").Value = 4146.7
$ws.Range("K132").Value = 12440.1
$ws.Range("M132").Value = -9910.099999999999
# Row 136
$ws.Range("H136").Value = 2041
$ws.Range("I136").Value = 1931.9231
$ws.Range("K136").Value = 5795.7693
$ws.Range("M136").Value = -3245.7693

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1233.3334
$ws.Range("I20").Value = 1233.3334
$ws.Range("K20").Value = 1233.3334
$ws.Range("M20").Value = -986.3334
# Row 86
$ws.Range("H86").Value = 2999.2778
$ws.Range("I86").Value = 1590.8
$ws.Range("J86").Value = 4759.875
$ws.Range("K86").Value = 1590.8
$ws.Range("L86").Value = 4759.875
$ws.Range("M86").Value = -467.8
$ws.Range("N86").Value = -7005.875
# Row 89
$ws.Range("H89").Value = 2999.2778
$ws.Range("I89").Value = 1590.8
$ws.Range("J89").Value = 4759.875
$ws.Range("K89").Value = 7954
$ws.Range("L89").Value = 23799.375
$ws.Range("M89").Value = -2338
$ws.Range("N89").Value = -35031.375
# Row 127
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()
$ws.Range("N127").ClearContents()
# Row 134
$ws.Range("H134").Value = 4879.4
$ws.Range("I134").Value = 1552.6
$ws.Range("J134").Value = 14859.8
$ws.Range("K134").Value = 4657.799999999999
$ws.Range("L134").Value = 44579.39999999999
$ws.Range("M134").Value = -2122.799999999999
$ws.Range("N134").Value = -49649.39999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3354.5854
$ws.Range("I31").Value = 1975.6842
$ws.Range("J31").Value = 4545.4546
$ws.Range("K31").Value = 1975.6842
$ws.Range("L31").Value = 4545.4546
$ws.Range("M31").Value = -1680.6842
$ws.Range("N31").Value = -5135.4546
# Row 34
$ws.Range("H34").Value = 3354.5854
$ws.Range("I34").Value = 1975.6842
$ws.Range("J34").Value = 4545.4546
$ws.Range("K34").Value = 1975.6842
$ws.Range("L34").Value = 4545.4546
$ws.Range("M34").Value = -1773.6842
$ws.Range("N34").Value = -4949.4546
# Row 93
$ws.Range("H93").Value = 7140.25
$ws.Range("I93").Value = 7140.25
$ws.Range("K93").Value = 7140.25
$ws.Range("M93").Value = -5268.25

$ws = $wb.Worksheets.Item("CUL")
# Row 70
$ws.Range("H70").Value = 1400
$ws.Range("I70").Value = 1400
$ws.Range("K70").Value = 1400
$ws.Range("M70").Value = -1130
# Row 73
$ws.Range("H73").Value = 1400
$ws.Range("I73").Value = 1400
$ws.Range("K73").Value = 1400
$ws.Range("M73").Value = -464
# Row 81
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("N81").ClearContents()
# Row 84
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
# Row 83
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
# Row 113
$ws.Range("H113").Value = 9229.625
$ws.Range("I113").Value = 9393
$ws.Range("K113").Value = 9393
$ws.Range("M113").Value = -7223
# Row 122
$ws.Range("H122").Value = 2309.4546
$ws.Range("I122").Value = 1197.1428
$ws.Range("K122").Value = 3591.4284
$ws.Range("M122").Value = -1141.4284
# Row 132
$ws.Range("H132").Value = 27401.268
$ws.Range("I132").Value = 32187.146
$ws.Range("K132").Value = 96561.43799999999
$ws.Range("M132").Value = -94031.43799999999

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3325
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 3325
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -3549
# Row 40
$ws.Range("H40").Value = 3845
$ws.Range("I40").Value = 3845
$ws.Range("K40").Value = 3845
$ws.Range("M40").Value = -3709
# Row 46
$ws.Range("H46").Value = 2108.85
$ws.Range("I46").Value = 538.2222
$ws.Range("K46").Value = 538.2222
$ws.Range("M46").Value = -350.2222
# Row 126
$ws.Range("H126").Value = 3325
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 3325
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -14915
# Row 132
$ws.Range("H132").Value = 2738.2354
$ws.Range("I132").Value = 1446.1538
$ws.Range("J132").Value = 6937.5
$ws.Range("K132").Value = 4338.4614
$ws.Range("L132").Value = 20812.5
$ws.Range("M132").Value = -1808.4614
$ws.Range("N132").Value = -25872.5
# Row 136
$ws.Range("H136").Value = 2012.7858
$ws.Range("I136").Value = 1860
$ws.Range("J136").Value = 2394.75
$ws.Range("K136").Value = 5580
$ws.Range("L136").Value = 7184.25
$ws.Range("M136").Value = -3030
$ws.Range("N136").Value = -12284.25

$ws = $wb.Worksheets.Item("WVR")
# Row 99
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
# Row 126
$ws.Range("H126").Value = 4165.6313
$ws.Range("I126").Value = 1949.75
$ws.Range("K126").Value = 5849.25
$ws.Range("M126").Value = -3379.25
# Row 132
$ws.Range("H132").Value = 2614.1428
$ws.Range("I132").Value = 2260
$ws.Range("K132").Value = 6780
$ws.Range("M132").Value = -4250
# Row 136
$ws.Range("H136").Value = 2200.4348
$ws.Range("I136").Value = 1925.3529
$ws.Range("J136").Value = 2979.8333
$ws.Range("K136").Value = 5776.0587
$ws.Range("L136").Value = 8939.499899999999
$ws.Range("M136").Value = -3226.0587
$ws.Range("N136").Value = -14039.4999
